$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clothing = @{
    2 = "Blouse,Trunks"
    3 = "Gauchos,Parka"
    4 = "Parka,Caftan"
    5 = "Jumpsuit,Dress"
    6 = "Trunks,Caftan"
    7 = "Caftan,Trunks"
    8 = "Trunks,Caftan"
    9 = "Halter,Blouse"
    10 = "Parka,Gauchos"
    11 = "Parka,Jumpsuit"
    12 = "Parka,Sweatpants"
    13 = "Parka,Jumpsuit"
    14 = "Caftan,Cutoffs"
    15 = "Caftan,Halter"
    16 = "Caftan,Parka"
    17 = "Kaftan,Trunks"
    18 = "Caftan,Trunks"
    19 = "Turtleneck,Jodhpurs"
    20 = "Jumpsuit,Blazer"
    21 = "Parka,Blouse"
    22 = "Caftan,Parka"
    23 = "Jumpsuit,Trunks"
    24 = "Trunks,Jumpsuit"
    25 = "Jumpsuit,Trunks"
    26 = "Kaftan,Blouse"
    27 = "Jumpsuit,Halter"
    28 = "Blazer,Halter"
    29 = "Trunks,Kaftan"
    30 = "Parka,Blouse"
    31 = "Turtleneck,Blazer"
    32 = "Caftan,Parka"
    33 = "Jumpsuit,Kaftan"
    34 = "Parka,Caftan"
    35 = "Jumpsuit,Kaftan"
    36 = "Parka,Gauchos"
    37 = "Trunks,Jumpsuit"
    38 = "Jumpsuit,Trunks"
    39 = "Jumpsuit,Parka"
    40 = "Jumpsuit,Dress"
    41 = "Parka,Trunks"
    42 = "Blouse,Parka"
    43 = "Caftan,Trunks"
    44 = "Caftan,Trunks"
    45 = "Caftan,Parka"
    46 = "Caftan,Sweatpants"
    47 = "Caftan,Parka"
    48 = "Blouse,Parka"
    49 = "Parka,Trunks"
    50 = "Jumpsuit,Tee"
    51 = "Caftan,Parka"
    52 = "Trunks,Kimono"
    53 = "Trunks,Sweatpants"
    54 = "Sweatpants,Trunks"
    55 = "Trunks,Jumpsuit"
    56 = "Trunks,Kimono"
    57 = "Jumpsuit,Halter"
    58 = "Jumpsuit,Kaftan"
    59 = "Caftan,Parka"
    60 = "Blouse,Tee"
    61 = "Blouse,Blazer"
    62 = "Jumpsuit,Parka"
    63 = "Caftan,Parka"
    64 = "Caftan,Parka"
    65 = "Caftan,Jumpsuit"
    66 = "Parka,Caftan"
    67 = "Jumpsuit,Blouse"
    68 = "Jumpsuit,Blouse"
    69 = "Jumpsuit,Dress"
    70 = "Jumpsuit,Tee"
    71 = "Blouse,Jumpsuit"
    72 = "Jumpsuit,Blazer"
    73 = "Cutoffs,Turtleneck"
    74 = "Trunks,Blouse"
    75 = "Caftan,Blazer"
    76 = "Caftan,Cutoffs"
    77 = "Blazer,Turtleneck"
    78 = "Caftan,Parka"
    79 = "Jodhpurs,Jumpsuit"
    80 = "Halter,Blazer"
    81 = "Jodhpurs,Cutoffs"
    82 = "Caftan,Sweatpants"
    83 = "Caftan,Blazer"
    84 = "Coat,Jumpsuit"
    85 = "Blouse,Trunks"
    86 = "Jumpsuit,Blouse"
    87 = "Jumpsuit,Sweatpants"
    88 = "Jumpsuit,Dress"
    89 = "Jodhpurs,Parka"
    90 = "Blouse,Parka"
    91 = "Blouse,Jumpsuit"
    92 = "Blazer,Top"
    93 = "Blazer,Blouse"
    94 = "Halter,Trunks"
    95 = "Jumpsuit,Trunks"
    96 = "Jodhpurs,Jumpsuit"
    97 = "Jumpsuit,Blouse"
    98 = "Parka,Sweatpants"
    99 = "Jumpsuit,Jodhpurs"
    100 = "Caftan,Halter"
    101 = "Parka,Caftan"
    102 = "Jumpsuit,Parka"
    103 = "Trunks,Caftan"
    104 = "Halter,Blazer"
    105 = "Caftan,Halter"
    106 = "Parka,Caftan"
    107 = "Dress,Trunks"
    108 = "Jumpsuit,Trunks"
    109 = "Jumpsuit,Trunks"
    110 = "Caftan,Trunks"
    111 = "Trunks,Jumpsuit"
    112 = "Jumpsuit,Parka"
    113 = "Jumpsuit,Trunks"
    114 = "Jumpsuit,Tee"
    115 = "Halter,Sweatpants"
    116 = "Jumpsuit,Dress"
    117 = "Caftan,Parka"
    118 = "Caftan,Jumpsuit"
    119 = "Parka,Gauchos"
    120 = "Gauchos,Cutoffs"
    121 = "Halter,Caftan"
    122 = "Blazer,Halter"
    123 = "Jumpsuit,Halter"
    124 = "Parka,Blouse"
    125 = "Jumpsuit,Caftan"
    126 = "Halter,Blouse"
    127 = "Blazer,Halter"
    128 = "Trunks,Jumpsuit"
    129 = "Halter,Jumpsuit"
    130 = "Sweatpants,Blouse"
    131 = "Blazer,Parka"
    132 = "Jumpsuit,Parka"
    133 = "Caftan,Trunks"
    134 = "Caftan,Trunks"
    135 = "Parka,Blouse"
}

foreach ($row in $clothing.Keys) {
    $ws.Range("G$row").Value = $clothing[$row]
}
